$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.73699999999999999
$ws.Range("C2").Value = 0.44500000000000001
$ws.Range("D2").Value = 0.70799999999999996

$ws.Range("B3").Value = 2.31
$ws.Range("C3").Value = 8.92
$ws.Range("D3").Value = 2.78

$ws.Range("B4").Value = 0.52400000000000002
$ws.Range("C4").Value = 1.1000000000000001
$ws.Range("D4").Value = 1.1499999999999999

$ws.Range("B5").Value = 0.17199999999999999
$ws.Range("C5").Value = 3.31
$ws.Range("D5").Value = 0.21

$ws.Range("B6").Value = 0.27700000000000002
$ws.Range("C6").Value = 2.92
$ws.Range("D6").Value = 0.32200000000000001

$ws.Range("B7").Value = 1.72
$ws.Range("C7").Value = 3.77
$ws.Range("D7").Value = 1.66

$ws.Range("B8").Value = 2.97
$ws.Range("C8").Value = 2.4300000000000002
$ws.Range("D8").Value = 0.76800000000000002

$wb.Save()
